# Update the "想去人数" (number of people interested) counts that changed
# between scrapes, as produced by the gh-pages data-refresh commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 297
$wsExhibition.Range("F3").Value = 95
$wsExhibition.Range("F4").Value = 1152
$wsExhibition.Range("F5").Value = 591

# Sheet "演出" (Performances)
$wsPerformance = $wb.Worksheets.Item("演出")
$wsPerformance.Range("F2").Value = 9

# Sheet "全部类型" (All types) - combined listing
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 297
$wsAll.Range("F3").Value = 95
$wsAll.Range("F4").Value = 1152
$wsAll.Range("F5").Value = 9
$wsAll.Range("F6").Value = 591
